$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 100-124 with revised values ---
# Row 100
$ws.Cells.Item(100, 4).Value = 44543
$ws.Cells.Item(100, 14).Value = 7000
$ws.Cells.Item(100, 15).Value = 7000
$ws.Cells.Item(100, 16).Value = 7000
$ws.Cells.Item(100, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(100, 19).Value = 700
# Row 101
$ws.Cells.Item(101, 4).Value = 44543
$ws.Cells.Item(101, 11).Value = 'Rainier'
$ws.Cells.Item(101, 12).Value = 'Primera'
$ws.Cells.Item(101, 13).Value = 150
# Row 102
$ws.Cells.Item(102, 4).Value = 44543
$ws.Cells.Item(102, 11).Value = 'Santina'
$ws.Cells.Item(102, 13).Value = 250
$ws.Cells.Item(102, 14).Value = 7000
$ws.Cells.Item(102, 15).Value = 8000
$ws.Cells.Item(102, 16).Value = 7400
$ws.Cells.Item(102, 19).Value = 740
# Row 103
$ws.Cells.Item(103, 4).Value = 44529
$ws.Cells.Item(103, 11).Value = 'Lapins'
$ws.Cells.Item(103, 13).Value = 200
$ws.Cells.Item(103, 14).Value = 15000
$ws.Cells.Item(103, 15).Value = 15000
$ws.Cells.Item(103, 16).Value = 15000
$ws.Cells.Item(103, 19).Value = 1500
# Row 104
$ws.Cells.Item(104, 4).Value = 44529
$ws.Cells.Item(104, 11).Value = 'Lapins'
$ws.Cells.Item(104, 12).Value = 'Segunda'
$ws.Cells.Item(104, 13).Value = 160
$ws.Cells.Item(104, 14).Value = 10000
$ws.Cells.Item(104, 15).Value = 10000
$ws.Cells.Item(104, 16).Value = 10000
$ws.Cells.Item(104, 19).Value = 1000
# Row 105
$ws.Cells.Item(105, 4).Value = 44529
$ws.Cells.Item(105, 11).Value = 'Rainier'
$ws.Cells.Item(105, 12).Value = 'Primera'
$ws.Cells.Item(105, 13).Value = 100
$ws.Cells.Item(105, 14).Value = 18000
$ws.Cells.Item(105, 15).Value = 18000
$ws.Cells.Item(105, 16).Value = 18000
$ws.Cells.Item(105, 19).Value = 1800
# Row 106
$ws.Cells.Item(106, 4).Value = 44175
$ws.Cells.Item(106, 11).Value = 'Rainier'
$ws.Cells.Item(106, 14).Value = 12000
$ws.Cells.Item(106, 15).Value = 12000
$ws.Cells.Item(106, 16).Value = 12000
$ws.Cells.Item(106, 19).Value = 1200
# Row 107
$ws.Cells.Item(107, 4).Value = 44175
$ws.Cells.Item(107, 11).Value = 'Santina'
$ws.Cells.Item(107, 14).Value = 9000
$ws.Cells.Item(107, 15).Value = 9000
$ws.Cells.Item(107, 16).Value = 9000
$ws.Cells.Item(107, 19).Value = 900
# Row 108
$ws.Cells.Item(108, 4).Value = 44175
$ws.Cells.Item(108, 11).Value = 'Santina'
$ws.Cells.Item(108, 12).Value = 'Segunda'
$ws.Cells.Item(108, 13).Value = 60
$ws.Cells.Item(108, 14).Value = 8000
$ws.Cells.Item(108, 15).Value = 8000
$ws.Cells.Item(108, 16).Value = 8000
$ws.Cells.Item(108, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(108, 19).Value = 800
# Row 109
$ws.Cells.Item(109, 4).Value = 44196
$ws.Cells.Item(109, 11).Value = 'Bing'
$ws.Cells.Item(109, 13).Value = 50
$ws.Cells.Item(109, 14).Value = 6000
$ws.Cells.Item(109, 15).Value = 6000
$ws.Cells.Item(109, 16).Value = 6000
$ws.Cells.Item(109, 19).Value = 600
# Row 110
$ws.Cells.Item(110, 4).Value = 44196
$ws.Cells.Item(110, 11).Value = 'Lapins'
$ws.Cells.Item(110, 13).Value = 100
$ws.Cells.Item(110, 18).Value = 'Provincia de Curicó'
# Row 111
$ws.Cells.Item(111, 4).Value = 44200
$ws.Cells.Item(111, 11).Value = 'Lapins'
$ws.Cells.Item(111, 12).Value = 'Primera'
$ws.Cells.Item(111, 13).Value = 250
$ws.Cells.Item(111, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(111, 18).Value = 'Provincia de Curicó'
# Row 112
$ws.Cells.Item(112, 4).Value = 44200
$ws.Cells.Item(112, 11).Value = 'Lapins'
$ws.Cells.Item(112, 13).Value = 130
$ws.Cells.Item(112, 17).Value = '$/caja 10 kilos'
# Row 113
$ws.Cells.Item(113, 4).Value = 44188
$ws.Cells.Item(113, 13).Value = 60
$ws.Cells.Item(113, 14).Value = 7000
$ws.Cells.Item(113, 16).Value = 7000
$ws.Cells.Item(113, 18).Value = 'Región del Maule'
$ws.Cells.Item(113, 19).Value = 700
# Row 114
$ws.Cells.Item(114, 4).Value = 44188
$ws.Cells.Item(114, 11).Value = 'Bing'
$ws.Cells.Item(114, 12).Value = 'Segunda'
$ws.Cells.Item(114, 13).Value = 80
$ws.Cells.Item(114, 14).Value = 6000
$ws.Cells.Item(114, 15).Value = 6000
$ws.Cells.Item(114, 16).Value = 6000
$ws.Cells.Item(114, 18).Value = 'Región del Maule'
$ws.Cells.Item(114, 19).Value = 600
# Row 115
$ws.Cells.Item(115, 4).Value = 44188
$ws.Cells.Item(115, 11).Value = 'Santina'
$ws.Cells.Item(115, 13).Value = 160
$ws.Cells.Item(115, 14).Value = 7000
$ws.Cells.Item(115, 15).Value = 7000
$ws.Cells.Item(115, 16).Value = 7000
$ws.Cells.Item(115, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(115, 19).Value = 700
# Row 116
$ws.Cells.Item(116, 4).Value = 44195
$ws.Cells.Item(116, 11).Value = 'Bing'
$ws.Cells.Item(116, 13).Value = 500
$ws.Cells.Item(116, 14).Value = 6000
$ws.Cells.Item(116, 15).Value = 7000
$ws.Cells.Item(116, 16).Value = 6600
$ws.Cells.Item(116, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(116, 19).Value = 660
$ws.Cells.Item(116, 20).Value = 10
# Row 117
$ws.Cells.Item(117, 4).Value = 44536
$ws.Cells.Item(117, 11).Value = 'Lapins'
$ws.Cells.Item(117, 13).Value = 180
$ws.Cells.Item(117, 14).Value = 10000
$ws.Cells.Item(117, 15).Value = 10000
$ws.Cells.Item(117, 16).Value = 10000
$ws.Cells.Item(117, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(117, 19).Value = 1000
$ws.Cells.Item(117, 20).Value = 10
# Row 118
$ws.Cells.Item(118, 4).Value = 44536
$ws.Cells.Item(118, 11).Value = 'Royal Dawn'
$ws.Cells.Item(118, 13).Value = 100
$ws.Cells.Item(118, 14).Value = 10000
$ws.Cells.Item(118, 15).Value = 10000
$ws.Cells.Item(118, 16).Value = 10000
$ws.Cells.Item(118, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(118, 19).Value = 1000
$ws.Cells.Item(118, 20).Value = 10
# Row 119
$ws.Cells.Item(119, 4).Value = 44511
$ws.Cells.Item(119, 11).Value = 'Early Burlat'
$ws.Cells.Item(119, 13).Value = 50
$ws.Cells.Item(119, 14).Value = 15000
$ws.Cells.Item(119, 15).Value = 15000
$ws.Cells.Item(119, 16).Value = 15000
$ws.Cells.Item(119, 17).Value = '$/bandeja 5 kilos'
$ws.Cells.Item(119, 19).Value = 3000
$ws.Cells.Item(119, 20).Value = 5
# Row 120
$ws.Cells.Item(120, 4).Value = 44511
$ws.Cells.Item(120, 11).Value = 'Early Burlat'
$ws.Cells.Item(120, 17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(120, 19).Value = 2857
$ws.Cells.Item(120, 20).Value = 7
# Row 121
$ws.Cells.Item(121, 4).Value = 44511
$ws.Cells.Item(121, 11).Value = 'Early Burlat'
$ws.Cells.Item(121, 12).Value = 'Primera'
$ws.Cells.Item(121, 13).Value = 20
$ws.Cells.Item(121, 14).Value = 3000
$ws.Cells.Item(121, 15).Value = 3000
$ws.Cells.Item(121, 16).Value = 3000
$ws.Cells.Item(121, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(121, 19).Value = 3000
$ws.Cells.Item(121, 20).Value = 1
# Row 122
$ws.Cells.Item(122, 4).Value = 44194
$ws.Cells.Item(122, 11).Value = 'Bing'
$ws.Cells.Item(122, 12).Value = 'Primera'
$ws.Cells.Item(122, 13).Value = 350
$ws.Cells.Item(122, 14).Value = 6000
$ws.Cells.Item(122, 15).Value = 6000
$ws.Cells.Item(122, 16).Value = 6000
$ws.Cells.Item(122, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(122, 19).Value = 600
# Row 123
$ws.Cells.Item(123, 4).Value = 44518
$ws.Cells.Item(123, 13).Value = 100
$ws.Cells.Item(123, 14).Value = 20000
$ws.Cells.Item(123, 15).Value = 20000
$ws.Cells.Item(123, 16).Value = 20000
$ws.Cells.Item(123, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(123, 19).Value = 2000
# Row 124
$ws.Cells.Item(124, 4).Value = 44518
$ws.Cells.Item(124, 11).Value = 'Royal Dawn'
$ws.Cells.Item(124, 12).Value = 'Segunda'
$ws.Cells.Item(124, 13).Value = 30
$ws.Cells.Item(124, 14).Value = 28000
$ws.Cells.Item(124, 15).Value = 28000
$ws.Cells.Item(124, 16).Value = 28000
$ws.Cells.Item(124, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(124, 19).Value = 1867
$ws.Cells.Item(124, 20).Value = 15

# --- Append new rows 125-127 ---
# Row 125
$ws.Cells.Item(125, 1).Value = 5
$ws.Cells.Item(125, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(125, 3).Value = 'Maule'
$ws.Cells.Item(125, 4).Value = 44518
$ws.Cells.Item(125, 5).Value = 7
$ws.Cells.Item(125, 6).Value = 'Fruta'
$ws.Cells.Item(125, 7).Value = 100103
$ws.Cells.Item(125, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(125, 9).Value = 100103001
$ws.Cells.Item(125, 10).Value = 'Cereza'
$ws.Cells.Item(125, 11).Value = 'Santina'
$ws.Cells.Item(125, 12).Value = 'Segunda'
$ws.Cells.Item(125, 13).Value = 150
$ws.Cells.Item(125, 14).Value = 15000
$ws.Cells.Item(125, 15).Value = 18000
$ws.Cells.Item(125, 16).Value = 17000
$ws.Cells.Item(125, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(125, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(125, 19).Value = 1700
$ws.Cells.Item(125, 20).Value = 10
$ws.Cells.Item(125, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 126
$ws.Cells.Item(126, 1).Value = 5
$ws.Cells.Item(126, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(126, 3).Value = 'Maule'
$ws.Cells.Item(126, 4).Value = 44540
$ws.Cells.Item(126, 5).Value = 7
$ws.Cells.Item(126, 6).Value = 'Fruta'
$ws.Cells.Item(126, 7).Value = 100103
$ws.Cells.Item(126, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(126, 9).Value = 100103001
$ws.Cells.Item(126, 10).Value = 'Cereza'
$ws.Cells.Item(126, 11).Value = 'Royal Dawn'
$ws.Cells.Item(126, 12).Value = 'Primera'
$ws.Cells.Item(126, 13).Value = 200
$ws.Cells.Item(126, 14).Value = 9000
$ws.Cells.Item(126, 15).Value = 9000
$ws.Cells.Item(126, 16).Value = 9000
$ws.Cells.Item(126, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(126, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(126, 19).Value = 900
$ws.Cells.Item(126, 20).Value = 10
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
# Row 127
$ws.Cells.Item(127, 1).Value = 5
$ws.Cells.Item(127, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(127, 3).Value = 'Maule'
$ws.Cells.Item(127, 4).Value = 44540
$ws.Cells.Item(127, 5).Value = 7
$ws.Cells.Item(127, 6).Value = 'Fruta'
$ws.Cells.Item(127, 7).Value = 100103
$ws.Cells.Item(127, 8).Value = 'Frutos de hueso (carozo)'
$ws.Cells.Item(127, 9).Value = 100103001
$ws.Cells.Item(127, 10).Value = 'Cereza'
$ws.Cells.Item(127, 11).Value = 'Santina'
$ws.Cells.Item(127, 12).Value = 'Primera'
$ws.Cells.Item(127, 13).Value = 300
$ws.Cells.Item(127, 14).Value = 8000
$ws.Cells.Item(127, 15).Value = 8000
$ws.Cells.Item(127, 16).Value = 8000
$ws.Cells.Item(127, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(127, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(127, 19).Value = 800
$ws.Cells.Item(127, 20).Value = 10
$ws.Cells.Item(127, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
